$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Cilantro" dataset for the
# Vega Modelo de Temuco market. It belongs chronologically before the
# existing row 575, so insert a fresh row there and push the rest of the
# table (old rows 575:681) down to 576:682.
$ws.Rows(575).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(575, 1).Value = 10
$ws.Cells.Item(575, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(575, 3).Value = "La Araucanía"
$ws.Cells.Item(575, 4).Value = 45258
$ws.Cells.Item(575, 5).Value = 9
$ws.Cells.Item(575, 6).Value = 100112040
$ws.Cells.Item(575, 7).Value = "Cilantro"
$ws.Cells.Item(575, 8).Value = "Sin especificar"
$ws.Cells.Item(575, 9).Value = "Primera"
$ws.Cells.Item(575, 10).Value = 75
$ws.Cells.Item(575, 11).Value = 6000
$ws.Cells.Item(575, 12).Value = 6000
$ws.Cells.Item(575, 13).Value = 6000
$ws.Cells.Item(575, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(575, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(575, 16).Value = 3000
$ws.Cells.Item(575, 17).Value = 2
$ws.Cells.Item(575, 18).Value = "Hortaliza"
